# Generate Report for Handback
# - updates the localization-status report with handback results:
#     * Overview "Status" column text changes from "Ready for handoff"
#       to "Handed back: in sync with en-US"
#     * zh-cn / de-de sheets get their "Latest Target File" / "Latest
#       Handback File" columns populated with file names (+ hyperlinks
#       for the handback file), and the "Latest Handback DateTime" gets
#       a real timestamp instead of the 0001-01-01 placeholder.
#     * a few columns are widened so the new, longer text fits.

$wb = $excel.ActiveWorkbook

$fileA = "5f91c3e9-9e6b-44c1-867a-aa6687cdeb9b"
$fileB = "f68abb70-02dd-4f90-9f99-fc68408bfc4e"

$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/82d5af4419224d768ae1247394beeb927e6153d1/e2e/$fileA.md"
$urlB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/82d5af4419224d768ae1247394beeb927e6153d1/e2e/$fileB.md"

# ---------------------------------------------------------------------
# Overview sheet: status text + column widths
#
# "Ready for handoff" is a single shared string reused by every status
# cell in the workbook (both locale columns on the Overview sheet, and
# the Status column on each locale sheet) - handback flips all of them
# to the new status text at once.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# widen the two status columns to fit the longer text
$overview.Columns.Item(5).ColumnWidth = 29.15
$overview.Columns.Item(6).ColumnWidth = 29.15

# ---------------------------------------------------------------------
# Helper: fill in one locale sheet's handback info
# ---------------------------------------------------------------------
function Update-LocaleSheet($sheetName, $targetFileA, $targetFileB, $handbackDateTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Status column - same shared "Handed back" text as the Overview sheet
    $ws.Range("C2").Value = "Handed back: in sync with en-US"
    $ws.Range("C3").Value = "Handed back: in sync with en-US"

    # Latest Target File / Latest Handback File / Latest Handback DateTime
    $ws.Range("I2").Value = "$fileA.md"
    $ws.Range("J2").Value = $targetFileA
    $ws.Range("K2").Value = $handbackDateTime

    $ws.Range("I3").Value = "$fileB.md"
    $ws.Range("J3").Value = $targetFileB
    $ws.Range("K3").Value = $handbackDateTime

    $ws.Range("I2").Style = "HyperLink"
    $ws.Range("I3").Style = "HyperLink"

    # rebuild hyperlinks in row order so relationship ids line up the
    # same way Excel would after a fresh report generation
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $urlA, "", "", "$fileA.md")
    $ws.Hyperlinks.Add($ws.Range("I2"), $urlA, "", "", "$fileA.md")
    $ws.Hyperlinks.Add($ws.Range("A3"), $urlB, "", "", "$fileB.md")
    $ws.Hyperlinks.Add($ws.Range("I3"), $urlB, "", "", "$fileB.md")

    # widen Status / Latest Target File / Latest Handback File columns
    $ws.Columns.Item(3).ColumnWidth = 29.15
    $ws.Columns.Item(9).ColumnWidth = 39.15
    $ws.Columns.Item(10).ColumnWidth = 39.15
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
Update-LocaleSheet "zh-cn" `
    "$fileA.a15e7ac57faef1c115c06271285e187902fffea4.zh-cn.xlf" `
    "$fileB.33527cb3cfe33ce0663b4d4a2ec45449e3561d5d.zh-cn.xlf" `
    "2016-08-28 18:46:36"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
Update-LocaleSheet "de-de" `
    "$fileA.a15e7ac57faef1c115c06271285e187902fffea4.de-de.xlf" `
    "$fileB.33527cb3cfe33ce0663b4d4a2ec45449e3561d5d.de-de.xlf" `
    "2016-08-28 18:46:42"

Write-Output "Handback report generated"
